# Append 6 new ScannedParameter / ScannedParameterUnit rows to Sheet1
# (rows 80-85), matching the onsite-changes commit that added the option
# to represent averaged data (NIhold, HfHoldTime, hw_SdCalibFreqStart,
# hw_SdRabiDepth, hw_XvNi, ZS_HF parameters with their units).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A80").Value = "NIhold"
$ws.Range("B80").Value = "ms"

$ws.Range("A81").Value = "HfHoldTime"
$ws.Range("B81").Value = "ms"

$ws.Range("A82").Value = "hw_SdCalibFreqStart"
$ws.Range("B82").Value = "Hz"

$ws.Range("A83").Value = "hw_SdRabiDepth"
$ws.Range("B83").Value = "V"

$ws.Range("A84").Value = "hw_XvNi"
$ws.Range("B84").Value = "V"

$ws.Range("A85").Value = "ZS_HF"
$ws.Range("B85").Value = "V"

# Match the author's final on-screen selection/scroll position.
$ws.Range("C76").Select()
